# Adjust offshore wind deployment for VCEA scenario
#
# The model previously assumed all 2,600 MW of Dominion's "2024-2026" phase
# and all 5,200 MW of the full VCEA mandate were still outstanding at the
# start of the schedule. In reality the model/dispatch engine already
# builds some offshore wind capacity in the start year (2019), so this
# edit carves that amount (568 MW) out of the remaining required buildout,
# and reworks the annual schedule so it is driven off the (now adjusted)
# per-phase totals rather than a FORECAST.LINEAR curve fit.

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$wind  = $wb.Worksheets.Item("Required offshore wind")
$pmccs = $wb.Worksheets.Item("PMCCS")

# ---------------------------------------------------------------------
# 1. Make room for the new "already built in 2019" explanatory block.
#    Inserting 3 rows above the old row 6 pushes "Required Offshore Wind"
#    (and everything below it) down to row 9, matching the new layout.
# ---------------------------------------------------------------------
$wind.Rows("4:6").Insert()

# New note explaining the model already builds some offshore wind in the
# start year, and the MW amount (568) that must be netted out.
$wind.Range("A4").Value2 = "Model Chooses to Build Offshore Wind Capacity in the start year, so we must adjust for that here (MW built 2019)"
$wind.Rows(4).RowHeight = 35.65

$wind.Range("A5").Value2 = 568
$wind.Range("A5").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. Net the already-built 568 MW out of the two phase targets.
#    B10 (2024-2026 phase, was a literal 2600) becomes 2600 - A5.
#    B11 (2026 cumulative target, was a literal 2600) becomes the
#    remainder needed to reach the final 5,200 MW (B12) total.
# ---------------------------------------------------------------------
$wind.Range("B10").Formula = "=2600-A5"
$wind.Range("B11").Formula = "=B12-B10-A5"

# ---------------------------------------------------------------------
# 3. Rebuild the per-year schedule (previously driven by
#    FORECAST.LINEAR off the two phase points) so that each year's
#    incremental build (column B) is simply the phase total spread
#    evenly across its years, and the cumulative total (column C) is
#    just a running sum of column B.
# ---------------------------------------------------------------------
$wind.Range("B24:B26").Formula = "=`$B`$10/3"
$wind.Range("B27:B35").Formula = "=`$B`$11/9"

$wind.Range("C24").Formula = "=B24"
$wind.Range("C25:C35").Formula = "=C24+B25"

# A couple of incidental number-format touches the author left behind
# while reworking this block.
$wind.Range("E24:E25").NumberFormat = "_(* #,##0_);_(* (#,##0);_(* ""-""??_);_(@_)"
$wind.Range("E33").NumberFormat = "0"

$wind.Rows(15).RowHeight = 28.5

# ---------------------------------------------------------------------
# 4. Force a full recalculation so the TRANSPOSE array spill on PMCCS
#    (row 14) picks up the reworked "Required offshore wind" numbers -
#    those spilled cells are cached literals in the OOXML, not live
#    formulas, so they need an explicit recalc to refresh.
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------
# 5. Restore the author's final selection/active-sheet state: "About"
#    ends up the active tab (rather than "PMCCS"), with "Required
#    offshore wind" scrolled down to the reworked schedule.
# ---------------------------------------------------------------------
$wind.Activate()
$wind.Range("E33").Select()

$about.Activate()
$about.Range("A1").Select()
